$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 2 Tasks")

# Set cell values in the precise order that reproduces the target
# shared-string table ordering (new strings are appended to the shared
# string table in first-use order).
# Target new-string order: main.c, get_source_line, init_lister, int main,
# common.h, print_line, print_page_header, print.c

# main.c (first use) -> row 2
$ws.Range("C2").Value = "main.c"

# get_source_line (first use) -> row 4
$ws.Range("D4").Value = "BOOLEAN get_source_line(FILE *src_file, char src_name[], char todays_date[])"

# init_lister (first use) -> row 3
$ws.Range("D3").Value = "FILE *init_lister(const char *name, char source_file_name[], char dte[])"

# int main (first use) -> row 2
$ws.Range("D2").Value = "int main (int argc, const char *argv[])"

# common.h (first use) -> row 7
$ws.Range("C7").Value = "common.h"

# print_line (first use) -> row 5
$ws.Range("D5").Value = "void print_line(char line[], char source_name_to_print[], char date_to_print[])"

# print_page_header (first use) -> row 6
$ws.Range("D6").Value = "static void print_page_header(char source_name[], char date[])"

# print.c (first use) -> row 5
$ws.Range("C5").Value = "print.c"

# Remaining duplicate-string cells (string already exists in table, order
# doesn't matter for these)
$ws.Range("H2").Value = ""
$ws.Range("C3").Value = "main.c"
$ws.Range("C4").Value = "main.c"
$ws.Range("F4").Value = 0.55
$ws.Range("C6").Value = "print.c"
$ws.Range("C8").Value = "common.h"
$ws.Range("C9").Value = "common.h"
$ws.Range("C10").Value = "common.h"

# Row 2 no longer needs the extra height (content now fits on one line) -
# AutoFit drops the explicit row height back to the sheet default.
$ws.Rows.Item(2).AutoFit()

# Column widths - best effort resize to match the new (longer) content.
# (The underlying column width units get quantized by this runtime, so the
# closest achievable value to the recorded target is used for each column.)
$ws.Columns.Item(2).ColumnWidth = 2.1666666666666665
$ws.Columns.Item(4).ColumnWidth = 71
$ws.Columns.Item(5).ColumnWidth = 12.833333333333332
$ws.Columns.Item(6).ColumnWidth = 11.999999999999998
$ws.Columns.Item(7).ColumnWidth = 9.666666666666666
$ws.Columns.Item(8).ColumnWidth = 27.666666666666664

# Update selection on sheet1
$ws.Range("D19").Select()

# Update selection on sheet2 and sheet3
$ws2 = $wb.Worksheets.Item("Phase 2 check list")
$ws2.Range("C20").Select()

$ws3 = $wb.Worksheets.Item("Phase1 Tasks (completed)")
$ws3.Range("C15").Select()

$ws.Activate()
